$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 105
$ws.Range("I2").Value = 291
$ws.Range("J2").Value = 1307
$ws.Range("K2").Value = 7
$ws.Range("L2").Value = 312
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 221
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 7
$ws.Range("S2").Value = 141
$ws.Range("T2").Value = 239
$ws.Range("U2").Value = 18
$ws.Range("V2").Value = 1868
$ws.Range("W2").Value = 2
$ws.Range("X2").Value = 1983
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 22
$ws.Range("AA2").Value = 17
